$wb = $excel.ActiveWorkbook

# --- Insert a new "Suggestions" sheet between Sheet1 and Sheet2 ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$newSheet = $wb.Worksheets.Add($sheet2)
$newSheet.Name = "Suggestions"

# --- Populate the Suggestions sheet. Write order chosen so that the
# shared-string table is built up in the same order as the target file. ---
$newSheet.Range("A2").Value = "RUN MOTOR"
$newSheet.Range("A1").Value = "STOP MOTOR "
$newSheet.Range("A5").Value = "MOTOR STOPPED  <%  <s"
$newSheet.Range("A3").Value = "FORWARD AT ###%"
$newSheet.Range("A4").Value = "REVERSE AT ###%"
$newSheet.Range("A6").Value = "For ##s"
$newSheet.Range("A8").Value = "STOPPED"
$newSheet.Range("A13").Value = "##s remaining at"
$newSheet.Range("A14").Value = "###% forward speed"
$newSheet.Range("A15").Value = "###% reverse speed"
$newSheet.Range("A12").Value = "0s No time set"
$newSheet.Range("A11").Value = "TRY THESE!"
$newSheet.Range("A16").Formula = '="  0% Increase speed"'
$newSheet.Range("A17").Formula = '="  0% No speed set"'

# Highlight the "try these" block with the yellow fill style.
$newSheet.Range("A12:A17").Interior.Color = 65535

# Column A is wide enough to show the suggestion text.
$newSheet.Columns.Item(1).ColumnWidth = 28.71

# --- View state ---
# Sheet1 selection changes from U10 to B4:H4.
$sheet1.Range("B4:H4").Select() | Out-Null

# Suggestions tab ends up active/selected, scrolled down a bit with A11 selected.
$newSheet.Activate()
$newSheet.Range("A11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
